# Update cryptocurrency price (D) and volume (E) figures to reflect the latest
# scrape from the GitHub Actions job that refreshes this workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "26.320.09" },
    @{ Cell = "E2"; Value = "  +0.10%  " },
    @{ Cell = "D3"; Value = "1.692.11" },
    @{ Cell = "E3"; Value = "  +0.65%  " },
    @{ Cell = "E4"; Value = "  +0.20%  " },
    @{ Cell = "D5"; Value = "217.97" },
    @{ Cell = "E5"; Value = "  -0.30%  " },
    @{ Cell = "D6"; Value = "0.5418" },
    @{ Cell = "E6"; Value = "  +2.69%  " },
    @{ Cell = "E7"; Value = "  +0.15%  " },
    @{ Cell = "D8"; Value = "0.2740" },
    @{ Cell = "E8"; Value = "  +1.29%  " },
    @{ Cell = "D9"; Value = "0.06457" },
    @{ Cell = "E9"; Value = "  -0.52%  " },
    @{ Cell = "E10"; Value = "  -1.45%  " },
    @{ Cell = "D11"; Value = "0.07657" },
    @{ Cell = "E11"; Value = "  +1.61%  " },
    @{ Cell = "D12"; Value = "1.712.34" },
    @{ Cell = "E12"; Value = "  +1.93%  " },
    @{ Cell = "D13"; Value = "4.541" },
    @{ Cell = "E13"; Value = "  +0.10%  " },
    @{ Cell = "D14"; Value = "0.5805" },
    @{ Cell = "E14"; Value = "  -0.16%  " },
    @{ Cell = "D15"; Value = "0.000008403" },
    @{ Cell = "E15"; Value = "  -1.40%  " },
    @{ Cell = "D16"; Value = "66.96" },
    @{ Cell = "E16"; Value = "  +3.52%  " },
    @{ Cell = "D17"; Value = "26.379.26" },
    @{ Cell = "E17"; Value = "  +0.22%  " },
    @{ Cell = "E19"; Value = "  +0.16%  " },
    @{ Cell = "D20"; Value = "10.89" },
    @{ Cell = "E20"; Value = "  +0.02%  " },
    @{ Cell = "D21"; Value = "190.72" },
    @{ Cell = "E21"; Value = "  +0.14%  " },
    @{ Cell = "D22"; Value = "6.284" },
    @{ Cell = "E22"; Value = "  +1.22%  " },
    @{ Cell = "E23"; Value = "  +0.15%  " },
    @{ Cell = "D24"; Value = "149.24" },
    @{ Cell = "E24"; Value = "  +2.59%  " },
    @{ Cell = "D25"; Value = "0.1286" },
    @{ Cell = "E25"; Value = "  +3.23%  " },
    @{ Cell = "D26"; Value = "7.857" },
    @{ Cell = "E26"; Value = "  +0.68%  " },
    @{ Cell = "E27"; Value = "  +0.36%  " },
    @{ Cell = "D28"; Value = "0.06332" },
    @{ Cell = "E28"; Value = "  -2.75%  " },
    @{ Cell = "D29"; Value = "1.382" },
    @{ Cell = "E29"; Value = "  +1.96%  " },
    @{ Cell = "D30"; Value = "1.325" },
    @{ Cell = "E30"; Value = "  -0.46%  " },
    @{ Cell = "D31"; Value = "3.607" },
    @{ Cell = "E31"; Value = "  +0.26%  " },
    @{ Cell = "D32"; Value = "3.590" },
    @{ Cell = "E32"; Value = "  -0.37%  " },
    @{ Cell = "D33"; Value = "1.684" },
    @{ Cell = "E33"; Value = "  +1.50%  " },
    @{ Cell = "E34"; Value = "  -0.12%  " },
    @{ Cell = "D35"; Value = "0.6197" },
    @{ Cell = "E35"; Value = "  -0.62%  " },
    @{ Cell = "D36"; Value = "2.417" },
    @{ Cell = "E36"; Value = "  +0.64%  " },
    @{ Cell = "D37"; Value = "2.773" },
    @{ Cell = "E37"; Value = "  +1.35%  " },
    @{ Cell = "D38"; Value = "0.01656" },
    @{ Cell = "E38"; Value = "  +1.86%  " },
    @{ Cell = "D39"; Value = "1.113.81" },
    @{ Cell = "E39"; Value = "  +0.01%  " },
    @{ Cell = "D40"; Value = "6.099" },
    @{ Cell = "D41"; Value = "0.8871" },
    @{ Cell = "E41"; Value = "  +1.17%  " },
    @{ Cell = "E42"; Value = "  -0.09%  " },
    @{ Cell = "D43"; Value = "101.08" },
    @{ Cell = "D44"; Value = "1.845.57" },
    @{ Cell = "E44"; Value = "  +0.83%  " },
    @{ Cell = "E45"; Value = "  -0.16%  " },
    @{ Cell = "D46"; Value = "57.73" },
    @{ Cell = "E46"; Value = "  +1.14%  " },
    @{ Cell = "D47"; Value = "8.191" },
    @{ Cell = "E47"; Value = "  +0.21%  " },
    @{ Cell = "D48"; Value = "1.004" },
    @{ Cell = "E48"; Value = "  -0.30%  " },
    @{ Cell = "D49"; Value = "0.05282" },
    @{ Cell = "E49"; Value = "  +0.16%  " },
    @{ Cell = "D50"; Value = "0.4301" },
    @{ Cell = "E50"; Value = "  +0.24%  " },
    @{ Cell = "D51"; Value = "6.070" },
    @{ Cell = "E51"; Value = "  -0.36%  " }
)

foreach ($u in $updates) {
    $cellRange = $ws.Range($u.Cell)
    if ($u.Cell.StartsWith("D")) {
        # Price column values are plain numeric-looking text (e.g. "217.97" or
        # thousand-grouped values like "26.320.09"). Force the cell to text format
        # first so Excel does not reinterpret the string as a floating point number
        # (which would introduce binary rounding artifacts), then restore the
        # default "Normal" style so no extra formatting is left behind.
        $cellRange.NumberFormat = "@"
        $cellRange.Value = $u.Value
        $cellRange.Style = "Normal"
    } else {
        $cellRange.Value = $u.Value
    }
}

